# Pojek_Deliverable1.docx edit script
# - Adds green (92D050) font color to 6 "feature title" list items
# - Merges the spell-check-split runs in the "Default roles available..."
#   and "Bees can be stinged" paragraphs back into plain runs
# - Moves the _GoBack bookmark from the trailing empty paragraph to just
#   after the "Encrypt Message" run

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Green color for the six list headings (pPr/rPr AND run rPr both gain
#    <w:color w:val="92D050"/> right before <w:sz w:val="28"/>)
# ---------------------------------------------------------------------------
$greenColor = 5296274  # 0x0050D092 -> RGB 92D050

$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Send Message" -or $t -eq "Create Hive" -or `
        $t -eq "Manage Hive (update and delete)" -or `
        $t -eq "Search Hive and Topic" -or `
        $t -eq "Reply to Conversation" -or `
        $t -eq "Encrypt Message") {
        $p.Range.Font.Color = $greenColor
    }
    $i++
}

# ---------------------------------------------------------------------------
# 2) "Default roles a" + <tab> + "vailable" + " are Queen..." -> one run
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Default roles a" + [char]9 + "vailable are Queen and Bee.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Default roles available are Queen and Bee.", 2)

# ---------------------------------------------------------------------------
# 3) "Bees can be " + "stinged" -> one run (keep the " " / "(kicked)" runs)
# ---------------------------------------------------------------------------
$i = 1
$beesIdx = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Bees can be stinged (kicked)")) {
        $beesIdx = $i
    }
    $i++
}

$p = $d.Paragraphs.Item($beesIdx)
$paraStart = $p.Range.Start
$rngAll = $d.Range($paraStart, $paraStart + 28)
$rngAll.Text = "Bees can be stingedPLACEHOLDER (kicked)"
$null = $d.Content.Find.Execute("stingedPLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "stinged", 2)

$splitA = $d.Range($paraStart + 19, $paraStart + 20)
$splitA.Font.Bold = 1
$splitA.Font.Bold = 0
$splitB = $d.Range($paraStart + 20, $paraStart + 28)
$splitB.Font.Bold = 1
$splitB.Font.Bold = 0

# ---------------------------------------------------------------------------
# 4) Move the _GoBack bookmark to right after "Encrypt Message"
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$i = 1
$encIdx = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Encrypt Message") {
        $encIdx = $i
    }
    $i++
}

$pEnc = $d.Paragraphs.Item($encIdx)
$insPos = $pEnc.Range.End - 1
$marker = $d.Range($insPos, $insPos)
$marker.InsertAfter("Z")
$rngMarker = $d.Range($insPos, $insPos + 1)
$d.Bookmarks.Add("_GoBack", $rngMarker)
$rngMarker2 = $d.Range($insPos, $insPos + 1)
$rngMarker2.Delete()
